# Update "想去人数" (want-to-go count) figures to match the refreshed scrape.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 199
$ws1.Range("F4").Value = 165
$ws1.Range("F6").Value = 18475
$ws1.Range("F7").Value = 383
$ws1.Range("F8").Value = 275
$ws1.Range("F10").Value = 6930
$ws1.Range("F11").Value = 112
$ws1.Range("F13").Value = 165
$ws1.Range("F14").Value = 20
$ws1.Range("F17").Value = 221
$ws1.Range("F18").Value = 167
$ws1.Range("F19").Value = 1307
$ws1.Range("F20").Value = 275
$ws1.Range("F24").Value = 36
$ws1.Range("F26").Value = 285
$ws1.Range("F27").Value = 1011
$ws1.Range("F28").Value = 4
$ws1.Range("F29").Value = 133
$ws1.Range("F30").Value = 5184
$ws1.Range("F31").Value = 544
$ws1.Range("F32").Value = 7
$ws1.Range("F33").Value = 56
$ws1.Range("F34").Value = 15
$ws1.Range("F36").Value = 12170
$ws1.Range("F37").Value = 1296
$ws1.Range("F38").Value = 14
$ws1.Range("F39").Value = 45
$ws1.Range("F40").Value = 215
$ws1.Range("F41").Value = 296
$ws1.Range("F42").Value = 3936

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 28

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 199
$ws4.Range("F4").Value = 165
$ws4.Range("F6").Value = 18475
$ws4.Range("F7").Value = 383
$ws4.Range("F8").Value = 275
$ws4.Range("F10").Value = 6930
$ws4.Range("F11").Value = 112
$ws4.Range("F13").Value = 165
$ws4.Range("F14").Value = 20
$ws4.Range("F17").Value = 221
$ws4.Range("F18").Value = 167
$ws4.Range("F19").Value = 1307
$ws4.Range("F20").Value = 275
$ws4.Range("F24").Value = 36
$ws4.Range("F26").Value = 285
$ws4.Range("F27").Value = 1011
$ws4.Range("F28").Value = 4
$ws4.Range("F29").Value = 133
$ws4.Range("F30").Value = 5184
$ws4.Range("F31").Value = 544
$ws4.Range("F33").Value = 7
$ws4.Range("F34").Value = 28
$ws4.Range("F35").Value = 56
$ws4.Range("F36").Value = 15
$ws4.Range("F38").Value = 12170
$ws4.Range("F39").Value = 1296
$ws4.Range("F40").Value = 14
$ws4.Range("F41").Value = 45
$ws4.Range("F42").Value = 215
$ws4.Range("F43").Value = 296
$ws4.Range("F44").Value = 3936

